$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140:190 down by one.
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with data (mostly copied from the row that used to be
# at 140, now at 141, with a few updated values per the weekly price update).
$ws.Cells.Item(140, 1).Value = 10
$ws.Cells.Item(140, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(140, 3).Value = "La Araucanía"
$ws.Cells.Item(140, 4).Value = 44559
$ws.Cells.Item(140, 4).NumberFormat = $ws.Cells.Item(141, 4).NumberFormat
$ws.Cells.Item(140, 5).Value = 9
$ws.Cells.Item(140, 6).Value = 100112039
$ws.Cells.Item(140, 7).Value = "Ciboulette"
$ws.Cells.Item(140, 8).Value = "Sin especificar"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 65
$ws.Cells.Item(140, 11).Value = 8000
$ws.Cells.Item(140, 12).Value = 8000
$ws.Cells.Item(140, 13).Value = 8000
$ws.Cells.Item(140, 14).Value = "`$/docena de atados"
$ws.Cells.Item(140, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(140, 16).Value = 2667
$ws.Cells.Item(140, 17).Value = 3
$ws.Cells.Item(140, 18).Value = "Hortaliza"
